# Append the "2021年" data row (row 12) to the bottom of the table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 1281732.7
$ws.Range("D12").Value = 1478226.1
$ws.Range("E12").Value = 5078234.1
$ws.Range("F12").Value = 2697909.8
$ws.Range("G12").Value = 6993975.1
$ws.Range("H12").Value = 11286100
$ws.Range("I12").Value = 1769182
$ws.Range("J12").Value = 593534.2
$ws.Range("K12").Value = 3437218.6
$ws.Range("M12").Value = 49034884.2
$ws.Range("N12").Value = 12712003.4
$ws.Range("O12").Value = 4243368.7
$ws.Range("P12").Value = 719569.6
$ws.Range("Q12").Value = 1449888.9
$ws.Range("T12").Value = 304525
$ws.Range("U12").Value = 75100171.8

# C12, L12, R12 and S12 stay blank for this year (same as the blank cells in
# the rows above), so they are intentionally left untouched.

# Match the formatting of the existing year-label cells (bold, centered,
# bordered) by copying the style from the cell directly above (A11).
$ws.Range("A11").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
